# Daily attendance processing - 2026-01-01 07:38:21
# Normalizes the "Recorded By" (column G) entries so that the literal
# "System" marker is listed last among the comma-separated recorder names,
# instead of first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $current = $cell.Value2

    if ($current -ne $null -and $current -like "System, *") {
        $parts = $current -split ", "
        if ($parts.Count -gt 1 -and $parts[0] -eq "System") {
            $rest = $parts[1..($parts.Count - 1)]
            $newValue = ($rest -join ", ") + ", System"
            $cell.Value = $newValue
        }
    }
}
